$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.7234226625123982
$ws.Cells.Item(2, 3).Value = 0.9983881055479173
$ws.Cells.Item(2, 4).Value = 0.8440790374608198
$ws.Cells.Item(2, 5).Value = 0.8441615919112424
$ws.Cells.Item(2, 6).Value = 0.8460098007712081
$ws.Cells.Item(2, 7).Value = 0.9134689692604095
$ws.Cells.Item(2, 8).Value = 0.9076387401784121
$ws.Cells.Item(2, 9).Value = 0.8356056580549259
$ws.Cells.Item(2, 10).Value = 0.9030659537338811
$ws.Cells.Item(2, 11).Value = 0.888450170348142
$ws.Cells.Item(2, 12).Value = 0.7739016458111628
$ws.Cells.Item(2, 13).Value = 0.9008466850943433
$ws.Cells.Item(3, 2).Value = 0.726499499692371
$ws.Cells.Item(3, 3).Value = 0.9982850073264816
$ws.Cells.Item(3, 4).Value = 0.8561735415842489
$ws.Cells.Item(3, 5).Value = 0.8517762485140665
$ws.Cells.Item(3, 6).Value = 0.7925152634961441
$ws.Cells.Item(3, 7).Value = 0.9155784520933978
$ws.Cells.Item(3, 8).Value = 0.9102836027641399
$ws.Cells.Item(3, 9).Value = 0.8513905397367761
$ws.Cells.Item(3, 10).Value = 0.8976464947223004
$ws.Cells.Item(3, 11).Value = 0.8887776815730692
$ws.Cells.Item(3, 12).Value = 0.8214125169718542
$ws.Cells.Item(3, 13).Value = 0.8968547316051689
$ws.Cells.Item(4, 2).Value = 0.7254668625566267
$ws.Cells.Item(4, 3).Value = 0.9983086411380399
$ws.Cells.Item(4, 4).Value = 0.860209005355893
$ws.Cells.Item(4, 5).Value = 0.8568068756303813
$ws.Cells.Item(4, 6).Value = 0.7900168701799485
$ws.Cells.Item(4, 7).Value = 0.9179993026800092
$ws.Cells.Item(4, 8).Value = 0.9019011823043968
$ws.Cells.Item(4, 9).Value = 0.8550366047826917
$ws.Cells.Item(4, 10).Value = 0.9042873950681861
$ws.Cells.Item(4, 11).Value = 0.8834359699690891
$ws.Cells.Item(4, 12).Value = 0.8207978606289676
$ws.Cells.Item(4, 13).Value = 0.8938877564122629
$ws.Cells.Item(5, 2).Value = 0.7348142402658966
$ws.Cells.Item(5, 3).Value = 0.9983487729484882
$ws.Cells.Item(5, 4).Value = 0.8461469332947871
$ws.Cells.Item(5, 5).Value = 0.8255504456441477
$ws.Cells.Item(5, 6).Value = 0.8567030848329048
$ws.Cells.Item(5, 7).Value = 0.9163921418507016
$ws.Cells.Item(5, 8).Value = 0.9047528037324314
$ws.Cells.Item(5, 9).Value = 0.8298749604837463
$ws.Cells.Item(5, 10).Value = 0.9021982590968226
$ws.Cells.Item(5, 11).Value = 0.8822075384693485
$ws.Cells.Item(5, 12).Value = 0.7940740706763054
$ws.Cells.Item(5, 13).Value = 0.895026763144732
$ws.Cells.Item(6, 2).Value = 0.7474127964897452
$ws.Cells.Item(6, 3).Value = 0.9983694953484776
$ws.Cells.Item(6, 4).Value = 0.8606963802075409
$ws.Cells.Item(6, 5).Value = 0.8597806836031372
$ws.Cells.Item(6, 6).Value = 0.7791862146529563
$ws.Cells.Item(6, 7).Value = 0.9125818182798482
$ws.Cells.Item(6, 8).Value = 0.9026144134368923
$ws.Cells.Item(6, 9).Value = 0.8267434475084821
$ws.Cells.Item(6, 10).Value = 0.8978705901813775
$ws.Cells.Item(6, 11).Value = 0.89123842253108
$ws.Cells.Item(6, 12).Value = 0.7693352996220322
$ws.Cells.Item(6, 13).Value = 0.8942585883854531
$ws.Cells.Item(7, 2).Value = 0.7322656272237416
$ws.Cells.Item(7, 3).Value = 0.9983203438708165
$ws.Cells.Item(7, 4).Value = 0.8464055985819057
$ws.Cells.Item(7, 5).Value = 0.8624182882202264
$ws.Cells.Item(7, 6).Value = 0.8132824550128536
$ws.Cells.Item(7, 7).Value = 0.9194849267454567
$ws.Cells.Item(7, 8).Value = 0.9158274328153214
$ws.Cells.Item(7, 9).Value = 0.817173067258826
$ws.Cells.Item(7, 10).Value = 0.8985091597854645
$ws.Cells.Item(7, 11).Value = 0.8864646555808759
$ws.Cells.Item(7, 12).Value = 0.8169768999302778
$ws.Cells.Item(7, 13).Value = 0.8957658718320376
$ws.Cells.Item(8, 2).Value = 0.7445765339300534
$ws.Cells.Item(8, 3).Value = 0.9983390111567574
$ws.Cells.Item(8, 4).Value = 0.8473541964638933
$ws.Cells.Item(8, 5).Value = 0.8536575896936978
$ws.Cells.Item(8, 6).Value = 0.7678470437017996
$ws.Cells.Item(8, 7).Value = 0.9171600025879917
$ws.Cells.Item(8, 8).Value = 0.9118551437125983
$ws.Cells.Item(8, 9).Value = 0.8306580789479501
$ws.Cells.Item(8, 10).Value = 0.9085610426343939
$ws.Cells.Item(8, 11).Value = 0.8855794234110242
$ws.Cells.Item(8, 12).Value = 0.8730275953176031
$ws.Cells.Item(8, 13).Value = 0.901912547535804
$ws.Cells.Item(9, 2).Value = 0.7468982569150833
$ws.Cells.Item(9, 3).Value = 0.9984297786939019
$ws.Cells.Item(9, 4).Value = 0.8520125490703265
$ws.Cells.Item(9, 5).Value = 0.8530729219643283
$ws.Cells.Item(9, 6).Value = 0.8110853149100256
$ws.Cells.Item(9, 7).Value = 0.9091691288532321
$ws.Cells.Item(9, 8).Value = 0.912573587252883
$ws.Cells.Item(9, 9).Value = 0.8172340832373253
$ws.Cells.Item(9, 10).Value = 0.9010565458227712
$ws.Cells.Item(9, 11).Value = 0.8824428600437716
$ws.Cells.Item(9, 12).Value = 0.7906693791053538
$ws.Cells.Item(9, 13).Value = 0.8932007111286004
$ws.Cells.Item(10, 2).Value = 0.7394505702631159
$ws.Cells.Item(10, 3).Value = 0.9983375839942238
$ws.Cells.Item(10, 4).Value = 0.8531242392197438
$ws.Cells.Item(10, 5).Value = 0.8488451453555409
$ws.Cells.Item(10, 6).Value = 0.7975120501285347
$ws.Cells.Item(10, 7).Value = 0.9124847685904073
$ws.Cells.Item(10, 8).Value = 0.9106119844548626
$ws.Cells.Item(10, 9).Value = 0.8398590002411812
$ws.Cells.Item(10, 10).Value = 0.9045375697993089
$ws.Cells.Item(10, 11).Value = 0.8885173295390446
$ws.Cells.Item(10, 12).Value = 0.8207818061722505
$ws.Cells.Item(10, 13).Value = 0.9026487815566511
$ws.Cells.Item(11, 2).Value = 0.7334274466356776
$ws.Cells.Item(11, 3).Value = 0.9982276924791275
$ws.Cells.Item(11, 4).Value = 0.8516288305285902
$ws.Cells.Item(11, 5).Value = 0.834042880292194
$ws.Cells.Item(11, 6).Value = 0.8325899742930591
$ws.Cells.Item(11, 7).Value = 0.9143956590033356
$ws.Cells.Item(11, 8).Value = 0.9069667739409147
$ws.Cells.Item(11, 9).Value = 0.8388471919677991
$ws.Cells.Item(11, 10).Value = 0.90574876570531
$ws.Cells.Item(11, 11).Value = 0.8797586076577694
$ws.Cells.Item(11, 12).Value = 0.8082788246302888
$ws.Cells.Item(11, 13).Value = 0.9078937702144948
$ws.Cells.Item(12, 2).Value = 0.7575016038308019
$ws.Cells.Item(12, 3).Value = 0.9983584776537174
$ws.Cells.Item(12, 4).Value = 0.8536477511335625
$ws.Cells.Item(12, 5).Value = 0.8554092635283028
$ws.Cells.Item(12, 6).Value = 0.8248015745501285
$ws.Cells.Item(12, 7).Value = 0.9149656012767426
$ws.Cells.Item(12, 8).Value = 0.9007179657362375
$ws.Cells.Item(12, 9).Value = 0.8300325450581382
$ws.Cells.Item(12, 10).Value = 0.9079694157199891
$ws.Cells.Item(12, 11).Value = 0.8847520080774351
$ws.Cells.Item(12, 12).Value = 0.8459322593666287
$ws.Cells.Item(12, 13).Value = 0.9147028968333315
$ws.Cells.Item(13, 2).Value = 0.7388752712280416
$ws.Cells.Item(13, 3).Value = 0.9982996214708268
$ws.Cells.Item(13, 4).Value = 0.8487930221234899
$ws.Cells.Item(13, 5).Value = 0.8490097530031191
$ws.Cells.Item(13, 6).Value = 0.8475698907455014
$ws.Cells.Item(13, 7).Value = 0.9072553179635381
$ws.Cells.Item(13, 8).Value = 0.9059790226648522
$ws.Cells.Item(13, 9).Value = 0.840396613469061
$ws.Cells.Item(13, 10).Value = 0.8982656911528677
$ws.Cells.Item(13, 11).Value = 0.8881286523318519
$ws.Cells.Item(13, 12).Value = 0.7854665883820777
$ws.Cells.Item(13, 13).Value = 0.8943116100112496
$ws.Cells.Item(14, 2).Value = 0.741111602538979
$ws.Cells.Item(14, 3).Value = 0.9984228712272387
$ws.Cells.Item(14, 4).Value = 0.8510225837619062
$ws.Cells.Item(14, 5).Value = 0.8461850236545848
$ws.Cells.Item(14, 6).Value = 0.8558973329048843
$ws.Cells.Item(14, 7).Value = 0.9168888475529101
$ws.Cells.Item(14, 8).Value = 0.9035530812062728
$ws.Cells.Item(14, 9).Value = 0.8373285186761742
$ws.Cells.Item(14, 10).Value = 0.9094836900232626
$ws.Cells.Item(14, 11).Value = 0.880853073317389
$ws.Cells.Item(14, 12).Value = 0.7621119408462075
$ws.Cells.Item(14, 13).Value = 0.8973619505314619
$ws.Cells.Item(15, 2).Value = 0.7309316560073043
$ws.Cells.Item(15, 3).Value = 0.9983462040559273
$ws.Cells.Item(15, 4).Value = 0.8477607384133166
$ws.Cells.Item(15, 5).Value = 0.8467557152629335
$ws.Cells.Item(15, 6).Value = 0.8070099614395887
$ws.Cells.Item(15, 7).Value = 0.917584594979296
$ws.Cells.Item(15, 8).Value = 0.9014709414781259
$ws.Cells.Item(15, 9).Value = 0.8107140214718591
$ws.Cells.Item(15, 10).Value = 0.9030870034320322
$ws.Cells.Item(15, 11).Value = 0.8822859789941562
$ws.Cells.Item(15, 12).Value = 0.7403661792227808
$ws.Cells.Item(15, 13).Value = 0.8952896354220233
$ws.Cells.Item(16, 2).Value = 0.7323113336646822
$ws.Cells.Item(16, 3).Value = 0.9983092690895549
$ws.Cells.Item(16, 4).Value = 0.847543916040291
$ws.Cells.Item(16, 5).Value = 0.8497621341847379
$ws.Cells.Item(16, 6).Value = 0.7838022172236504
$ws.Cells.Item(16, 7).Value = 0.9134826730360018
$ws.Cells.Item(16, 8).Value = 0.9130057827320749
$ws.Cells.Item(16, 9).Value = 0.8336589120322548
$ws.Cells.Item(16, 10).Value = 0.9017256654312544
$ws.Cells.Item(16, 11).Value = 0.8865083707497574
$ws.Cells.Item(16, 12).Value = 0.8340152654948442
$ws.Cells.Item(16, 13).Value = 0.9017225001421362
$ws.Cells.Item(17, 2).Value = 0.7381549158235151
$ws.Cells.Item(17, 3).Value = 0.998285578191495
$ws.Cells.Item(17, 4).Value = 0.8514795274033049
$ws.Cells.Item(17, 5).Value = 0.8358589995737282
$ws.Cells.Item(17, 6).Value = 0.7869770244215939
$ws.Cells.Item(17, 7).Value = 0.9129178977599494
$ws.Cells.Item(17, 8).Value = 0.9098262139806333
$ws.Cells.Item(17, 9).Value = 0.8194071168652969
$ws.Cells.Item(17, 10).Value = 0.8942425908788117
$ws.Cells.Item(17, 11).Value = 0.8781624751810654
$ws.Cells.Item(17, 12).Value = 0.834649416535173
$ws.Cells.Item(17, 13).Value = 0.8884607055453594
$ws.Cells.Item(18, 2).Value = 0.727185643689008
$ws.Cells.Item(18, 3).Value = 0.9983457473639166
$ws.Cells.Item(18, 4).Value = 0.8458469005812361
$ws.Cells.Item(18, 5).Value = 0.8334318132231187
$ws.Cells.Item(18, 6).Value = 0.8195276349614397
$ws.Cells.Item(18, 7).Value = 0.9101737728605935
$ws.Cells.Item(18, 8).Value = 0.8965523831779963
$ws.Cells.Item(18, 9).Value = 0.8177366242885873
$ws.Cells.Item(18, 10).Value = 0.8995882828066066
$ws.Cells.Item(18, 11).Value = 0.8860676936440965
$ws.Cells.Item(18, 12).Value = 0.7630018164471027
$ws.Cells.Item(18, 13).Value = 0.8925254839181576
$ws.Cells.Item(19, 2).Value = 0.7334882060960898
$ws.Cells.Item(19, 3).Value = 0.9982369975788473
$ws.Cells.Item(19, 4).Value = 0.8500706574662975
$ws.Cells.Item(19, 5).Value = 0.8342812508007152
$ws.Cells.Item(19, 6).Value = 0.8010178341902314
$ws.Cells.Item(19, 7).Value = 0.9143383727283184
$ws.Cells.Item(19, 8).Value = 0.9056115479157101
$ws.Cells.Item(19, 9).Value = 0.8209627840966416
$ws.Cells.Item(19, 10).Value = 0.8984622170957823
$ws.Cells.Item(19, 11).Value = 0.8888651119108323
$ws.Cells.Item(19, 12).Value = 0.7981610766577374
$ws.Cells.Item(19, 13).Value = 0.8989420588616703
$ws.Cells.Item(20, 2).Value = 0.7460884044675172
$ws.Cells.Item(20, 3).Value = 0.9984209302861928
$ws.Cells.Item(20, 4).Value = 0.8559305673899151
$ws.Cells.Item(20, 5).Value = 0.8552524203924029
$ws.Cells.Item(20, 6).Value = 0.8307189910025707
$ws.Cells.Item(20, 7).Value = 0.917615372311364
$ws.Cells.Item(20, 8).Value = 0.9092511116545529
$ws.Cells.Item(20, 9).Value = 0.8392646950017826
$ws.Cells.Item(20, 10).Value = 0.8989212122837843
$ws.Cells.Item(20, 11).Value = 0.8913945984973263
$ws.Cells.Item(20, 12).Value = 0.8020038255476863
$ws.Cells.Item(20, 13).Value = 0.8948942090802408
$ws.Cells.Item(21, 2).Value = 0.717840455509843
$ws.Cells.Item(21, 3).Value = 0.9983384973782453
$ws.Cells.Item(21, 4).Value = 0.8497240269620524
$ws.Cells.Item(21, 5).Value = 0.8627327509431941
$ws.Cells.Item(21, 6).Value = 0.8289725257069409
$ws.Cells.Item(21, 7).Value = 0.9227104810070162
$ws.Cells.Item(21, 8).Value = 0.9093136605480239
$ws.Cells.Item(21, 9).Value = 0.8289491511091936
$ws.Cells.Item(21, 10).Value = 0.9049684366296044
$ws.Cells.Item(21, 11).Value = 0.8847819741206201
$ws.Cells.Item(21, 12).Value = 0.8082788246302888
$ws.Cells.Item(21, 13).Value = 0.8957658718320376
